$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "WARNING: not found -> $old"
    }
}

# ---------------------------------------------------------------
# Title
Replace-Text "Unveiling the Enigmatic Universe" "Exploring the Marvels of the Microscopic World: Unveiling the Wonders of Chemistry"

# Author name
Replace-Text "Harper Cassandra" "Emily Carter"

# ---------------------------------------------------------------
# Email: "harpercassie98@gmail" + "." + "com"  ->  "emily" + "." + "carter@validedu" + "." + "org"
Replace-Text "harpercassie98@gmail" "emily"

$rEmail = $d.Content
$rEmail.Find.Execute("emily", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterEmily = $rEmail.End
$comRange = $d.Range($afterEmily + 1, $afterEmily + 4)
if ($comRange.Text -eq "com") {
    $comRange.Text = "carter@validedu"
} else {
    Write-Host "WARNING: email 'com' run not where expected, got: [$($comRange.Text)]"
}

$rEmail2 = $d.Content
$rEmail2.Find.Execute("carter@validedu", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertPt = $d.Range($rEmail2.End, $rEmail2.End)
$insertPt.InsertAfter(".org")

# ---------------------------------------------------------------
# Body paragraph (intro)
Replace-Text "In the ocean of human consciousness, the allure of understanding the unseen, the mystical, and the enigmatic has forever stoked the fires of our curiosity" "Within the boundless realms of science, chemistry emerges as a captivating force that unravels the intricacies of matter, transforming the ordinary into the extraordinary"

Replace-Text " Across eons, we have embarked on a fearless quest to unlock the secrets of the universe, peering through telescopes and microscopes, conducting unyielding experiments, and engaging in intellectual discourse" " From the smallest atoms to the most complex molecules, chemistry orchestrates a symphony of reactions that shape our world"

Replace-Text " Our relentless pursuit of knowledge has illuminated the cosmos, revealing planets, stars, and galaxies beyond our earthly abode. It has unveiled the mind-boggling world of the infinitely small, where atoms and subatomic particles perform their intricate dance" " As we embark on this journey into the microscopic realm, we will decipher the enigmatic language of chemical equations, unravel the mysteries of chemical bonds, and witness the awe-inspiring transformations that occur when substances interact"

Replace-Text "Yet, amidst the vast expanse of our newfound knowledge, unfathomable mysteries endure, taunting us with their elusiveness" "In the vast expanse of chemistry, we will delve into the wonders of the periodic table, where elements dance in harmonious arrangements, each possessing unique properties that contribute to the tapestry of our universe"

Replace-Text " The nature of dark matter and dark energy, the true depths of the black holes that swallow stardust, the genesis of the universe itself--these questions and countless others dance at the fringes of our understanding, beckoning us to unravel their enigmatic nature" " We will explore the fundamental principles that govern chemical reactions, examining how atoms rearrange themselves to form new substances with distinct characteristics"

Replace-Text " We stand at the threshold of profound discoveries, poised to push the boundaries of human knowledge further than ever before" " Through experimentation and observation, we will uncover the secrets that lie hidden within the molecular structures of matter, gaining insights into the forces that drive chemical change"

Replace-Text "As we venture into this uncharted territory, we must embrace the unknown with audacity and unyielding hope" "With unwavering curiosity, we will investigate the practical applications of chemistry, witnessing its transformative impact on fields ranging from medicine to engineering"

Replace-Text " We must be willing to challenge long-held beliefs, to entertain unconventional ideas, and to venture beyond the confines of established paradigms" " We will learn how chemistry enables the development of life-saving drugs, fuels our technological advancements, and shapes the materials that surround us"

Replace-Text " It is through this relentless pursuit of understanding, through our collective hunger for knowledge, that we will ultimately illuminate the darkness that shrouds the universe's most profound enigmas" " Furthermore, we will explore the intricate relationship between chemistry and the environment, examining how human activities can disrupt delicate ecosystems and the measures we can take to mitigate these effects"

Write-Host "DONE STEP1-2"

